$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest Kaspa buy-log entry as row 13, mirroring the existing
# rows (e.g. rows 10/12) that store the date as literal text rather than
# a serial date value.
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "04/22/2025"
$ws.Range("A13").Style = "Normal"

$ws.Range("B13").Value = 540.8340000000007
$ws.Range("C13").Value = 0.09244980899869448
$ws.Range("D13").Value = 50
